$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 = "Save", copying formatting from the G1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H13 with 0
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# H14 = 1
$ws.Cells.Item(14, 8).Value = 1
